$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 12: new Schottky diode part (order matches shared-string insertion order)
$ws.Range("B12").Value = "Schottky Diode"
$ws.Range("C12").Value = "0603/SOD-523F"
$ws.Range("F12").Value = "DIODE SCHOTTKY 20V 500MA 0603"
$ws.Range("G12").Value = "Comchip Technology"
$ws.Range("H12").Value = "CDBU0530"
$ws.Range("I12").Value = "641-1285-1-ND"
$ws.Range("E12").Value = "20V"
$ws.Range("K12").Value = 0.4
$ws.Range("L12").Value = 0.31
$ws.Range("M12").Value = 0.216
$ws.Range("N12").Value = 0.11

# Row 15: trailing space
$ws.Range("F15").Value = " "

# Row heights explicitly set to match default (forces customHeight flag)
$ws.Rows.Item(12).RowHeight = 14
$ws.Rows.Item(13).RowHeight = 14

# Unhide column A and D, adjust D width (closest reachable to the
# original bestFit value of 5.6640625 chars given this host's width
# quantization)
$ws.Columns.Item(1).Hidden = $false
$ws.Columns.Item(4).Hidden = $false
$ws.Columns.Item(4).ColumnWidth = 4.8

# View settings: zoom to 125% and move selection to F15
$ws.Range("F15").Select()
$excel.ActiveWindow.Zoom = 125
